$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 0.1636683333333333
$ws.Range("H2").Value = 0.491005
$ws.Range("I2").Value = 0.008639493057305454
$ws.Range("J2").Value = 0.008639493057305455
$ws.Range("M2").Value = 7.487621999999999
$ws.Range("N2").Value = 22.462866
$ws.Range("O2").Value = 0.1384395179233961
$ws.Range("P2").Value = 0.1384395179233961
$ws.Range("Q2").Value = 1.22548661337
$ws.Range("R2").Value = 11.02937952033
$ws.Range("S2").Value = 0.001196047253955894
$ws.Range("T2").Value = 0.001196047253955895
$ws.Range("G3").Value = 0.1636683333333333
$ws.Range("H3").Value = 0.491005
$ws.Range("I3").Value = 0.008639493057305454
$ws.Range("J3").Value = 0.008639493057305455
$ws.Range("O3").Value = 0.5916411627275552
$ws.Range("P3").Value = 0.5916411627275552
$ws.Range("Q3").Value = 5.23729304838
$ws.Range("R3").Value = 47.13563743542
$ws.Range("S3").Value = 0.005111479717800839
$ws.Range("T3").Value = 0.00511147971780084
$ws.Range("G4").Value = 0.1636683333333333
$ws.Range("H4").Value = 0.491005
$ws.Range("I4").Value = 0.008639493057305454
$ws.Range("J4").Value = 0.008639493057305455
$ws.Range("M4").Value = 14.59882166666667
$ws.Range("N4").Value = 43.796465
$ws.Range("O4").Value = 0.2699193193490487
$ws.Range("P4").Value = 0.2699193193490487
$ws.Range("Q4").Value = 2.389364810813889
$ws.Range("R4").Value = 21.504283297325
$ws.Range("S4").Value = 0.00233196608554872
$ws.Range("T4").Value = 0.00233196608554872
$ws.Range("I5").Value = 0.808839719627903
$ws.Range("J5").Value = 0.8088397196279031
$ws.Range("M5").Value = 7.487621999999999
$ws.Range("N5").Value = 22.462866
$ws.Range("O5").Value = 0.1384395179233961
$ws.Range("P5").Value = 0.1384395179233961
$ws.Range("Q5").Value = 114.731529059772
$ws.Range("R5").Value = 1032.583761537948
$ws.Range("S5").Value = 0.1119753808625817
$ws.Range("T5").Value = 0.1119753808625817
$ws.Range("I6").Value = 0.808839719627903
$ws.Range("J6").Value = 0.8088397196279031
$ws.Range("O6").Value = 0.5916411627275552
$ws.Range("P6").Value = 0.5916411627275552
$ws.Range("R6").Value = 4412.895006091752
$ws.Range("S6").Value = 0.4785428721808823
$ws.Range("T6").Value = 0.4785428721808823
$ws.Range("I7").Value = 0.808839719627903
$ws.Range("J7").Value = 0.8088397196279031
$ws.Range("M7").Value = 14.59882166666667
$ws.Range("N7").Value = 43.796465
$ws.Range("O7").Value = 0.2699193193490487
$ws.Range("P7").Value = 0.2699193193490487
$ws.Range("Q7").Value = 223.6952042033633
$ws.Range("R7").Value = 2013.25683783027
$ws.Range("S7").Value = 0.218321466584439
$ws.Range("T7").Value = 0.218321466584439
$ws.Range("G8").Value = 3.457711333333334
$ws.Range("H8").Value = 10.373134
$ws.Range("I8").Value = 0.1825207873147914
$ws.Range("J8").Value = 0.1825207873147914
$ws.Range("M8").Value = 7.487621999999999
$ws.Range("N8").Value = 22.462866
$ws.Range("O8").Value = 0.1384395179233961
$ws.Range("P8").Value = 0.1384395179233961
$ws.Range("Q8").Value = 25.890035449116
$ws.Range("R8").Value = 233.010319042044
$ws.Range("S8").Value = 0.02526808980685843
$ws.Range("T8").Value = 0.02526808980685843
$ws.Range("G9").Value = 3.457711333333334
$ws.Range("H9").Value = 10.373134
$ws.Range("I9").Value = 0.1825207873147914
$ws.Range("J9").Value = 0.1825207873147914
$ws.Range("O9").Value = 0.5916411627275552
$ws.Range("P9").Value = 0.5916411627275552
$ws.Range("Q9").Value = 110.644784855784
$ws.Range("R9").Value = 995.803063702056
$ws.Range("S9").Value = 0.107986810828872
$ws.Range("T9").Value = 0.107986810828872
$ws.Range("G10").Value = 3.457711333333334
$ws.Range("H10").Value = 10.373134
$ws.Range("I10").Value = 0.1825207873147914
$ws.Range("J10").Value = 0.1825207873147914
$ws.Range("M10").Value = 14.59882166666667
$ws.Range("N10").Value = 43.796465
$ws.Range("O10").Value = 0.2699193193490487
$ws.Range("P10").Value = 0.2699193193490487
$ws.Range("Q10").Value = 50.47851113014556
$ws.Range("R10").Value = 454.30660017131
$ws.Range("S10").Value = 0.04926588667906099
$ws.Range("T10").Value = 0.04926588667906099
